$wb = $excel.ActiveWorkbook

$wsCentenario = $wb.Worksheets.Item("CENTENARIO")
$wsCentenario.Range("G10").Value = 1.267
$wsCentenario.Range("G17").Value = 0.3

$wsElPenon = $wb.Worksheets.Item("EL_PEÑON")
$wsElPenon.Range("E10").Value = 2.286
$wsElPenon.Range("F10").Value = 0.286
$wsElPenon.Range("H10").Value = 0.2
$wsElPenon.Range("I10").Value = 0.2
$wsElPenon.Range("G14").Value = 3.9
$wsElPenon.Range("F15").Value = 2.125
$wsElPenon.Range("I15").Value = 0.2
$wsElPenon.Range("D18").Value = 0.667
$wsElPenon.Range("G18").Value = 0.8
$wsElPenon.Range("D24").Value = 3
$wsElPenon.Range("G24").Value = 1.1
$wsElPenon.Range("F30").Value = 4.25
$wsElPenon.Range("D32").Value = 4.143
$wsElPenon.Range("G32").Value = 1.3
$wsElPenon.Range("D40").Value = 9
$wsElPenon.Range("E40").Value = 10.6
$wsElPenon.Range("G40").Value = 4.4
$wsElPenon.Range("H40").Value = 0.2
$wsElPenon.Range("D43").Value = 0.25
$wsElPenon.Range("G43").Value = 0.5
$wsElPenon.Range("E46").Value = 4
$wsElPenon.Range("H46").Value = 0.2
$wsElPenon.Range("D53").Value = 6.857
$wsElPenon.Range("G53").Value = 4.1

$wsGranada = $wb.Worksheets.Item("GRANADA")
$wsGranada.Range("D10").Value = 9.4
$wsGranada.Range("G10").Value = 0.4
$wsGranada.Range("I16").Value = 2.6
$wsGranada.Range("D17").Value = 9.199999999999999
$wsGranada.Range("G17").Value = 0.1
$wsGranada.Range("D23").Value = 1.833
$wsGranada.Range("G23").Value = 0.633
$wsGranada.Range("I25").Value = 0.733
$wsGranada.Range("D29").Value = 2.417
$wsGranada.Range("G29").Value = 0.4
$wsGranada.Range("H30").Value = 0.4
$wsGranada.Range("D31").Value = 3.667
$wsGranada.Range("G31").Value = 0.5
$wsGranada.Range("H31").Value = 0.4
$wsGranada.Range("I31").Value = 0.7
$wsGranada.Range("D32").Value = 2
$wsGranada.Range("G32").Value = 0.3
$wsGranada.Range("D38").Value = 5.833
$wsGranada.Range("G38").Value = 0.1
$wsGranada.Range("I40").Value = 3
$wsGranada.Range("G42").Value = 0.533
$wsGranada.Range("D47").Value = 5.333
$wsGranada.Range("G47").Value = 2.4
$wsGranada.Range("G52").Value = 2.6
$wsGranada.Range("D68").Value = 9.182
$wsGranada.Range("G68").Value = 0.45
$wsGranada.Range("H68").Value = 0.9
$wsGranada.Range("I68").Value = 0.7
$wsGranada.Range("D69").Value = 4.5
$wsGranada.Range("G69").Value = 1.6
$wsGranada.Range("D73").Value = 3.5
$wsGranada.Range("G73").Value = 0.767
$wsGranada.Range("D75").Value = 1.077
$wsGranada.Range("G75").Value = 0.7
$wsGranada.Range("I75").Value = 1.667
$wsGranada.Range("I78").Value = 1.067
$wsGranada.Range("D79").Value = 3.125
$wsGranada.Range("G79").Value = 0.2
$wsGranada.Range("D81").Value = 4
$wsGranada.Range("G81").Value = 0.067
$wsGranada.Range("D86").Value = 0.545
$wsGranada.Range("G86").Value = 1.35
$wsGranada.Range("G91").Value = 1
$wsGranada.Range("D92").Value = 1.333
$wsGranada.Range("G92").Value = 1.067

$wsSanAntonio = $wb.Worksheets.Item("SAN_ANTONIO")
$wsSanAntonio.Range("D7").Value = 1.933
$wsSanAntonio.Range("G7").Value = 0.367
$wsSanAntonio.Range("G20").Value = 2.2
$wsSanAntonio.Range("D30").Value = 1.857
$wsSanAntonio.Range("G30").Value = 0.7
$wsSanAntonio.Range("D58").Value = 1.5
$wsSanAntonio.Range("G58").Value = 1.2
$wsSanAntonio.Range("D59").Value = 2.154
$wsSanAntonio.Range("G59").Value = 0.4
$wsSanAntonio.Range("D60").Value = 2.333
$wsSanAntonio.Range("G60").Value = 1.5
$wsSanAntonio.Range("D78").Value = 2

$wsSanFernando = $wb.Worksheets.Item("SAN_FERNANDO_PARQUE_DEL_PERRO")
$wsSanFernando.Range("D6").Value = 2.333
$wsSanFernando.Range("G6").Value = 0.7
$wsSanFernando.Range("E8").Value = 14
$wsSanFernando.Range("G18").Value = 1.2
$wsSanFernando.Range("E19").Value = 0
$wsSanFernando.Range("H19").Value = 0.1
$wsSanFernando.Range("D22").Value = 7.857
$wsSanFernando.Range("G22").Value = 0.5
$wsSanFernando.Range("E31").Value = 0.833
$wsSanFernando.Range("H31").Value = 0.2
$wsSanFernando.Range("D36").Value = 2.875
$wsSanFernando.Range("G36").Value = 0.1
$wsSanFernando.Range("G46").Value = 0.5
$wsSanFernando.Range("D48").Value = 2.3
$wsSanFernando.Range("E48").Value = 1.5
$wsSanFernando.Range("G48").Value = 0.7
$wsSanFernando.Range("H48").Value = 0.1
$wsSanFernando.Range("G55").Value = 0.7
$wsSanFernando.Range("G60").Value = 2.4
$wsSanFernando.Range("E65").Value = 3.143
$wsSanFernando.Range("H65").Value = 0.2
$wsSanFernando.Range("D67").Value = 6.875
$wsSanFernando.Range("G67").Value = 0.5
$wsSanFernando.Range("E72").Value = 3.375
$wsSanFernando.Range("H72").Value = 0.2
$wsSanFernando.Range("G77").Value = 0
$wsSanFernando.Range("D78").Value = 2.875
$wsSanFernando.Range("G78").Value = 0.25
$wsSanFernando.Range("D93").Value = 3.8
$wsSanFernando.Range("G93").Value = 1.2
$wsSanFernando.Range("D94").Value = 0.5
$wsSanFernando.Range("G94").Value = 0.65
$wsSanFernando.Range("D99").Value = 4.833
$wsSanFernando.Range("G99").Value = 0.2
